$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Número expediente" column into the header row (row 7) ---
# The original layout was:
#   A7=Folio Pase a Caja  B7=Nombre del Solicitante  C7=Folio de pago
#   D7=Fecha de Solicitud E7=Motivo de Cancelación
# The new layout inserts "Número expediente" as column C, shifting the
# following headers one column to the right:
#   A7=Folio Pase a Caja  B7=Nombre del Solicitante  C7=Número expediente
#   D7=Folio de pago      E7=Fecha de Solicitud      F7=Motivo de Cancelación
#
# Shift the existing header cells (value + formatting) one column to the
# right, starting from the rightmost column so values are not overwritten
# before they are copied.
$ws.Range("E7").Copy($ws.Range("F7"))
$ws.Range("D7").Copy($ws.Range("E7"))
$ws.Range("C7").Copy($ws.Range("D7"))

# Put the new header text in the now-free C7 cell (keeps the header style
# that was already there, since it was just copied from the old C7).
$ws.Range("C7").Value2 = "Número expediente"

# --- Column widths: nudge column C (now bestfit-ish) and E a bit ---
$ws.Columns("C:C").ColumnWidth = 17
$ws.Columns("E:E").ColumnWidth = 31.333333333333336

# --- Update the active selection as recorded in the saved file ---
[void]$ws.Range("D11").Select()
